# Remove pre existing sensors and some buildings
#
# The canonical edit deletes three previously-placed sensor/building
# picture icons from the site-plan slide (slide 5): "Picture 14",
# "Picture 25" and "Picture 28". All other content is left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

$s.Shapes.Item("Picture 14").Delete()
$s.Shapes.Item("Picture 25").Delete()
$s.Shapes.Item("Picture 28").Delete()
